# Fruta / hortaliza, semanal
# Insert a new week of data (3 rows) at the top of the Chirimoya / Femacal de
# La Calera block (rows 117-119), pushing the existing rows 117-134 down to
# 120-137.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 117:134 down by inserting 3 new rows before row 117.
$ws.Rows("117:119").Insert()

# Fill the 3 newly inserted rows with the new week's data
# (date 44522 = 2021-11-22), mirroring the layout used by every other row
# in this block.

# Row 117 - Especial
$ws.Cells.Item(117, 1).Value = 3
$ws.Cells.Item(117, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(117, 3).Value = "Coquimbo"
$ws.Cells.Item(117, 4).Value = 44522
$ws.Cells.Item(117, 5).Value = 5
$ws.Cells.Item(117, 6).Value = "Fruta"
$ws.Cells.Item(117, 7).Value = 100107
$ws.Cells.Item(117, 8).Value = "Otros"
$ws.Cells.Item(117, 9).Value = 100107002
$ws.Cells.Item(117, 10).Value = "Chirimoya"
$ws.Cells.Item(117, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(117, 12).Value = "Especial"
$ws.Cells.Item(117, 13).Value = 45
$ws.Cells.Item(117, 14).Value = 26000
$ws.Cells.Item(117, 15).Value = 26000
$ws.Cells.Item(117, 16).Value = 26000
$ws.Cells.Item(117, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(117, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(117, 19).Value = 2600
$ws.Cells.Item(117, 20).Value = 10

# Row 118 - Primera
$ws.Cells.Item(118, 1).Value = 3
$ws.Cells.Item(118, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(118, 3).Value = "Coquimbo"
$ws.Cells.Item(118, 4).Value = 44522
$ws.Cells.Item(118, 5).Value = 5
$ws.Cells.Item(118, 6).Value = "Fruta"
$ws.Cells.Item(118, 7).Value = 100107
$ws.Cells.Item(118, 8).Value = "Otros"
$ws.Cells.Item(118, 9).Value = 100107002
$ws.Cells.Item(118, 10).Value = "Chirimoya"
$ws.Cells.Item(118, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(118, 12).Value = "Primera"
$ws.Cells.Item(118, 13).Value = 55
$ws.Cells.Item(118, 14).Value = 23000
$ws.Cells.Item(118, 15).Value = 23000
$ws.Cells.Item(118, 16).Value = 23000
$ws.Cells.Item(118, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(118, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(118, 19).Value = 2300
$ws.Cells.Item(118, 20).Value = 10

# Row 119 - Segunda
$ws.Cells.Item(119, 1).Value = 3
$ws.Cells.Item(119, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(119, 3).Value = "Coquimbo"
$ws.Cells.Item(119, 4).Value = 44522
$ws.Cells.Item(119, 5).Value = 5
$ws.Cells.Item(119, 6).Value = "Fruta"
$ws.Cells.Item(119, 7).Value = 100107
$ws.Cells.Item(119, 8).Value = "Otros"
$ws.Cells.Item(119, 9).Value = 100107002
$ws.Cells.Item(119, 10).Value = "Chirimoya"
$ws.Cells.Item(119, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(119, 12).Value = "Segunda"
$ws.Cells.Item(119, 13).Value = 45
$ws.Cells.Item(119, 14).Value = 20000
$ws.Cells.Item(119, 15).Value = 20000
$ws.Cells.Item(119, 16).Value = 20000
$ws.Cells.Item(119, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(119, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(119, 19).Value = 2000
$ws.Cells.Item(119, 20).Value = 10
